$wb = $excel.ActiveWorkbook

# --- Rename the existing (only) sheet to Test_Content ---
$wb.Worksheets.Item(1).Name = "Test_Content"

# --- Add the new Test_Title sheet. Add() inserts it at the front and makes it active. ---
$titleWs = $wb.Worksheets.Add()
$titleWs.Name = "Test_Title"

# Re-fetch the content sheet by name now that sheet order/indices have changed.
$contentWs = $wb.Worksheets.Item("Test_Content")

# === Test_Content sheet: insert new column B ("id_test") for rows 2-10 ===
$contentWs.Range("B2").Value = "id_test"
$contentWs.Range("B3").Value = "id_test"
$contentWs.Range("B4").Value = "id_test"
$contentWs.Range("B5").Value = "id_test"
$contentWs.Range("B6").Value = "id_test"
$contentWs.Range("B7").Value = "id_test"
$contentWs.Range("B8").Value = "id_test"
$contentWs.Range("B9").Value = "id_test"
$contentWs.Range("B10").Value = "id_test"

$contentWs.Columns.Item(2).ColumnWidth = 17.42578125
$contentWs.Columns.Item(3).ColumnWidth = 92.85546875

# Select on the content sheet first so the LAST .Select() call in the
# script (on Test_Title, below) is what ends up as the truly active /
# tabSelected sheet+selection in the saved file.
$contentWs.Range("A2").Select()

# === Test_Title sheet content ===
$titleWs.Range("A1").Value = "DialogueTalkerNameID"
$titleWs.Range("B1").Value = "DialogueText"
$titleWs.Range("A2").Value = "id_000"
$titleWs.Range("B2").Value = "Long Dialogue Test"

$titleWs.Columns.Item(1).ColumnWidth = 25.42578125
$titleWs.Columns.Item(2).ColumnWidth = 40.5703125

$titleWs.Range("A1").Validation.Add(0)
$titleWs.Range("A1").Validation.IgnoreBlank = $false

$titleWs.Range("B2").Select()

Write-Host "Sheets:"
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Host $i $wb.Worksheets.Item($i).Name
}
